# Updates coin price/volume figures and re-ranks a handful of exchange tokens
# (rows 6-17) to match the refreshed coinranking.com snapshot referenced in the
# commit message. Values are written as text via NumberFormat "@" so that
# numeric-looking strings (prices, percentages) are preserved exactly as text,
# matching how the source sheet stores every data cell as a string.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $text
}

Set-TextValue "D2" '308.85'
Set-TextValue "E2" '0.34%'
Set-TextValue "D3" '41.30'
Set-TextValue "E3" '3.65%'
Set-TextValue "D4" '5.132'
Set-TextValue "E4" '0.90%'
Set-TextValue "D5" '0.07650'
Set-TextValue "E5" '-0.57%'
$ws.Range("B6").Value = 'GateToken'
$ws.Range("C6").Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
Set-TextValue "D6" '4.270'
Set-TextValue "E6" '0.67%'
$ws.Range("B7").Value = 'FTXToken'
$ws.Range("C7").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
Set-TextValue "D7" '1.616'
Set-TextValue "E7" '-0.63%'
$ws.Range("B8").Value = 'BTSEToken'
$ws.Range("C8").Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
Set-TextValue "D8" '2.504'
Set-TextValue "E8" '3.40%'
$ws.Range("B9").Value = 'MXToken'
$ws.Range("C9").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
Set-TextValue "D9" '0.9083'
Set-TextValue "E9" '-1.12%'
$ws.Range("B10").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C10").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
Set-TextValue "D10" '0.1156'
Set-TextValue "E10" '11.61%'
$ws.Range("B11").Value = 'WazirX'
$ws.Range("C11").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
Set-TextValue "D11" '0.1800'
Set-TextValue "E11" '0.78%'
$ws.Range("B12").Value = 'MandalaExchangeToken'
$ws.Range("C12").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
Set-TextValue "D12" '0.09168'
Set-TextValue "E12" '-1.25%'
$ws.Range("B13").Value = 'BitrueCoin'
$ws.Range("C13").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
Set-TextValue "D13" '0.04251'
Set-TextValue "E13" '-4.11%'
$ws.Range("B14").Value = 'BitMartToken'
$ws.Range("C14").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
Set-TextValue "D14" '0.1043'
Set-TextValue "E14" '-1.22%'
$ws.Range("B15").Value = 'BitForexToken'
$ws.Range("C15").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
Set-TextValue "D15" '0.001250'
Set-TextValue "E15" '-1.50%'
$ws.Range("B16").Value = 'TigerCash'
$ws.Range("C16").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
Set-TextValue "D16" '0.005752'
Set-TextValue "E16" '-1.47%'
$ws.Range("B17").Value = 'LEO'
$ws.Range("C17").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
Set-TextValue "D17" '3.354'
Set-TextValue "E17" '-0.31%'
Set-TextValue "E18" '0.66%'
Set-TextValue "D19" '6.687'
Set-TextValue "E19" '-3.11%'
Set-TextValue "D20" '0.1359'
Set-TextValue "E20" '0.98%'
Set-TextValue "D21" '0.2734'
Set-TextValue "E21" '1.09%'
Set-TextValue "D22" '0.04056'
Set-TextValue "D23" '0.001273'
Set-TextValue "E23" '5.78%'
Set-TextValue "D24" '0.004054'
Set-TextValue "E24" '-1.11%'
Set-TextValue "D25" '0.0001271'
Set-TextValue "E25" '-2.19%'
Set-TextValue "D38" '0.02436'
Set-TextValue "E38" '-1.86%'
Set-TextValue "D39" '0.05265'
Set-TextValue "E39" '1.44%'
Set-TextValue "D40" '0.007789'
Set-TextValue "E40" '-1.89%'
Set-TextValue "D41" '0.1303'
Set-TextValue "E41" '-1.22%'
Set-TextValue "D42" '0.006746'
Set-TextValue "E42" '-4.27%'
Set-TextValue "D43" '0.001950'
Set-TextValue "E43" '0.12%'
Set-TextValue "D44" '0.007570'
Set-TextValue "E44" '-5.06%'
Set-TextValue "E45" '0.34%'
Set-TextValue "D46" '0.00006904'
Set-TextValue "E46" '8.01%'
Set-TextValue "E47" '0.11%'
Set-TextValue "D48" '0.07956'
Set-TextValue "E48" '1,669.00%'
Set-TextValue "D49" '0.003001'
Set-TextValue "E49" '0.13%'
Set-TextValue "D50" '0.00002101'
Set-TextValue "E50" '0.11%'
Set-TextValue "D51" '0.0002001'
Set-TextValue "E51" '0.11%'
